$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$giordani = "198273 - Domingos Savio Giordani"
$media = "A média para a primeira avaliação será calculada a partir das notas das duas provas, P1 e P2, segundo a fórmula: M1=(P1+2xP2)/3. Alunos com nota final igual ou superior a 5,0 estão aprovados; inferior a 5,0 e igual ou superior a 3,0 estão de recuperação;"
$recuperacao = "A recuperação consistirá de uma prova envolvendo o assunto do semestre todo, à qual será atribuída nota NR. A média da segunda avaliação será calculada segunda a fórmula: M2=(M1+NR)/2. Alunos com nota M2 igual ou superior a 5,0 estarão aprovados, inferior a 5,0 estarão reprovados."

# The standalone row 13 (only B/C = "198273 - Domingos Savio Giordani", column A empty) is removed,
# shifting everything below it up by one row.
$ws.Rows(13).Delete()

# Row 10 ("Objetivos:") now shows the docente's name instead of the long objectives paragraph.
$ws.Range("B10").Value = $giordani
$ws.Range("C10").Value = $giordani

# Row 13 ("Programa resumido:", shifted up from old row 14) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:", shifted up from old row 16) now just shows the activation date, reused
# from row 8 ("Ativação:") so that it keeps being stored as text (not auto-converted to a date
# serial number) and keeps the existing cell style intact.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 ("Método:", shifted up from old row 19) now shows the docente's name.
$ws.Range("B18").Value = $giordani
$ws.Range("C18").Value = $giordani

# Row 19 ("Critério:", shifted up from old row 20) now shows "Duas provas escritas".
$ws.Range("B19").Value = "Duas provas escritas"
$ws.Range("C19").Value = "Duas provas escritas"

# Row 20 ("Norma de recuperação:", shifted up from old row 21) now holds the average-grade formula text.
$ws.Range("B20").Value = $media
$ws.Range("C20").Value = $media

# Row 21 ("Bibliografia:", shifted up from old row 22) now holds the recovery-exam rule text.
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao
